$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/D may contain values that look numeric (e.g. "5.980", "0.06020")
# -- a bare COM .Value assignment would auto-coerce those to Number and silently
# drop significant trailing zeros / reformat them. We force literal text by
# prefixing with an apostrophe (the standard Excel "treat as text" marker) and
# restore the original cell Style afterwards so no stray quotePrefix formatting
# style is left behind on the cell.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '27.566.61'
$ws.Range("E2").Value = '  -0.08%  '

Set-TextValue $ws.Range("D3") '1.752.10'
$ws.Range("E3").Value = '  -0.27%  '

$ws.Range("E4").Value = '  +0.04%  '

Set-TextValue $ws.Range("D5") '324.36'
$ws.Range("E5").Value = '  -0.17%  '

$ws.Range("E6").Value = '  +0.07%  '

Set-TextValue $ws.Range("D7") '0.4492'
$ws.Range("E7").Value = '  +0.57%  '

$ws.Range("E8").Value = '  -1.95%  '

Set-TextValue $ws.Range("D9") '0.07461'
$ws.Range("E9").Value = '  -0.58%  '

Set-TextValue $ws.Range("D10") '41.39'
$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("E11").Value = '  -2.35%  '

$ws.Range("E12").Value = '  +0.05%  '

Set-TextValue $ws.Range("D13") '20.72'
$ws.Range("E13").Value = '  -0.26%  '

Set-TextValue $ws.Range("D14") '5.980'
$ws.Range("E14").Value = '  -1.32%  '

Set-TextValue $ws.Range("D15") '7.146'
$ws.Range("E15").Value = '  -1.21%  '

Set-TextValue $ws.Range("D16") '1.752.82'
$ws.Range("E16").Value = '  -0.23%  '

Set-TextValue $ws.Range("D17") '93.62'
$ws.Range("E17").Value = '  +0.67%  '

Set-TextValue $ws.Range("D18") '0.00001056'
$ws.Range("E18").Value = '  -0.86%  '

Set-TextValue $ws.Range("D19") '0.06368'
$ws.Range("E19").Value = '  -0.89%  '

$ws.Range("E20").Value = '  +0.03%  '

Set-TextValue $ws.Range("D21") '17.14'
$ws.Range("E21").Value = '  +0.19%  '

$ws.Range("E22").Value = '  -2.14%  '

Set-TextValue $ws.Range("D23") '27.616.15'
$ws.Range("E23").Value = '  -0.05%  '

Set-TextValue $ws.Range("D24") '11.19'
$ws.Range("E24").Value = '  -0.53%  '

Set-TextValue $ws.Range("D25") '2.085'
$ws.Range("E25").Value = '  -0.07%  '

Set-TextValue $ws.Range("D26") '165.37'
$ws.Range("E26").Value = '  +1.49%  '

$ws.Range("E27").Value = '  -1.64%  '

Set-TextValue $ws.Range("D28") '1.956.27'
$ws.Range("E28").Value = '  -0.03%  '

Set-TextValue $ws.Range("D29") '2.094'
$ws.Range("E29").Value = '  -2.03%  '

Set-TextValue $ws.Range("D30") '125.35'
$ws.Range("E30").Value = '  -0.43%  '

Set-TextValue $ws.Range("D31") '1.089'
$ws.Range("E31").Value = '  -0.68%  '

Set-TextValue $ws.Range("D32") '0.09188'
$ws.Range("E32").Value = '  +1.03%  '

Set-TextValue $ws.Range("D33") '3.652'
$ws.Range("E33").Value = '  +0.35%  '

Set-TextValue $ws.Range("D34") '5.496'
$ws.Range("E34").Value = '  -1.44%  '

Set-TextValue $ws.Range("D35") '11.77'
$ws.Range("E35").Value = '  -3.55%  '

Set-TextValue $ws.Range("D36") '0.02286'
$ws.Range("E36").Value = '  -0.83%  '

Set-TextValue $ws.Range("B37") 'Hedera'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D37") '0.06020'
$ws.Range("E37").Value = '  +0.08%  '

Set-TextValue $ws.Range("B38") 'Algorand'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D38") '0.2087'
$ws.Range("E38").Value = '  -1.02%  '

Set-TextValue $ws.Range("D39") '0.6280'
$ws.Range("E39").Value = '  -2.15%  '

Set-TextValue $ws.Range("D40") '4.926'
$ws.Range("E40").Value = '  -0.36%  '

Set-TextValue $ws.Range("D41") '1.181'
$ws.Range("E41").Value = '  -0.91%  '

Set-TextValue $ws.Range("D42") '1.401'
$ws.Range("E42").Value = '  +0.24%  '

Set-TextValue $ws.Range("D43") '7.785'
$ws.Range("E43").Value = '  -1.00%  '

Set-TextValue $ws.Range("D44") '13.22'
$ws.Range("E44").Value = '  -0.67%  '

Set-TextValue $ws.Range("D45") '3.716'
$ws.Range("E45").Value = '  +0.14%  '

Set-TextValue $ws.Range("D46") '0.5862'
$ws.Range("E46").Value = '  -1.04%  '

Set-TextValue $ws.Range("D47") '121.96'
$ws.Range("E47").Value = '  -0.09%  '

Set-TextValue $ws.Range("D48") '1.932'
$ws.Range("E48").Value = '  -2.26%  '

Set-TextValue $ws.Range("D49") '0.06878'
$ws.Range("E49").Value = '  +0.07%  '

Set-TextValue $ws.Range("D50") '1.129'
$ws.Range("E50").Value = '  -3.27%  '

Set-TextValue $ws.Range("D51") '71.57'
$ws.Range("E51").Value = '  -1.89%  '

